# edit.ps1 - applies DNB Mastercard Demo -> Sheet data refresh
# (feat: increase merchant diversity in test data, PLAN 1.4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet/tab from "DNB Mastercard Demo" to "Sheet"
$ws.Name = "Sheet"

# 2. The transaction data rows (2..39): date serial, description, Inn, Ut
$transactions = @(
    @{ r = 2; date = 45762; name = "GITHUB.COM"; inn = $null; out = 129 },
    @{ r = 3; date = 45761; name = "JUST EAT"; inn = $null; out = 267 },
    @{ r = 4; date = 45759; name = "MENY CC VEST"; inn = $null; out = 534.6 },
    @{ r = 5; date = 45757; name = "POWER LAMBERTSETER"; inn = $null; out = 1899 },
    @{ r = 6; date = 45755; name = "JACK & JONES OSLO CITY"; inn = $null; out = 1199 },
    @{ r = 7; date = 45753; name = "G-SPORT STORO"; inn = $null; out = 1599 },
    @{ r = 8; date = 45752; name = "Innbetaling"; inn = 15000; out = $null },
    @{ r = 9; date = 45750; name = "COOP PRIX SAGENE"; inn = $null; out = 534.2 },
    @{ r = 10; date = 45748; name = "LINDEX OSLO CITY"; inn = $null; out = 699 },
    @{ r = 11; date = 45747; name = "SPOTIFY"; inn = $null; out = 129 },
    @{ r = 12; date = 45746; name = "Kiwi Torshov"; inn = $null; out = 312 },
    @{ r = 13; date = 45744; name = "netflix.com"; inn = $null; out = 179 },
    @{ r = 14; date = 45743; name = "kiwi grünerløkka"; inn = $null; out = 267.5 },
    @{ r = 15; date = 45741; name = "NILLE STORO"; inn = $null; out = 149 },
    @{ r = 16; date = 45740; name = "VINMONOPOLET GRÜNERLØKKA"; inn = $null; out = 567 },
    @{ r = 17; date = 45738; name = "Starbucks Aker Brygge"; inn = $null; out = 89 },
    @{ r = 18; date = 45736; name = "NORMAL MAJORSTUEN"; inn = $null; out = 199 },
    @{ r = 19; date = 45734; name = "GITHUB.COM"; inn = $null; out = 129 },
    @{ r = 20; date = 45733; name = "FLYING TIGER OSLO"; inn = $null; out = 89 },
    @{ r = 21; date = 45731; name = "REMA 1000 GRÜNERLØKKA"; inn = $null; out = 534.2 },
    @{ r = 22; date = 45730; name = "starbucks bogstadveien"; inn = $null; out = 75 },
    @{ r = 23; date = 45728; name = "MENY MAJORSTUEN"; inn = $null; out = 623.45 },
    @{ r = 24; date = 45726; name = "KICKS OSLO CITY"; inn = $null; out = 456 },
    @{ r = 25; date = 45724; name = "INTERSPORT CC VEST"; inn = $null; out = 1299 },
    @{ r = 26; date = 45722; name = "VITA KARL JOHAN"; inn = $null; out = 289 },
    @{ r = 27; date = 45721; name = "Innbetaling"; inn = 15000; out = $null },
    @{ r = 28; date = 45719; name = "COOP MEGA TORSHOV"; inn = $null; out = 756.4 },
    @{ r = 29; date = 45717; name = "ESSO HOVINBYEN"; inn = $null; out = 678 },
    @{ r = 30; date = 45716; name = "SPOTIFY FAMILY"; inn = $null; out = 169 },
    @{ r = 31; date = 45715; name = "KIWI MINIPRIS LØKKA"; inn = $null; out = 289.5 },
    @{ r = 32; date = 45714; name = "NETFLIX PREMIUM"; inn = $null; out = 219 },
    @{ r = 33; date = 45713; name = "REMA TORSHOV"; inn = $null; out = 345 },
    @{ r = 34; date = 45712; name = "JULA OSLO"; inn = $null; out = 567 },
    @{ r = 35; date = 45710; name = "VINMONOPOLET AKER BRYGGE"; inn = $null; out = 489 },
    @{ r = 36; date = 45708; name = "STARBUCKS AKER BRYGGE"; inn = $null; out = 95 },
    @{ r = 37; date = 45706; name = "BILTEMA ALF"; inn = $null; out = 789 },
    @{ r = 38; date = 45705; name = "GITHUB.COM"; inn = $null; out = 129 },
    @{ r = 39; date = 45703; name = "COOP OBS BYGG ALNA"; inn = $null; out = 1456 }

)

foreach ($t in $transactions) {
    $ws.Cells.Item($t.r, 1).Value = $t.date
    $ws.Cells.Item($t.r, 2).Value = $t.name

    if ($t.inn -ne $null) {
        $ws.Cells.Item($t.r, 5).Value = $t.inn
    } else {
        $ws.Cells.Item($t.r, 5).ClearContents()
    }

    if ($t.out -ne $null) {
        $ws.Cells.Item($t.r, 6).Value = $t.out
    } else {
        $ws.Cells.Item($t.r, 6).ClearContents()
    }
}

# 3. Apply the updated date/time number format (yyyy-mm-dd -> yyyy-mm-dd h:mm:ss)
#    to the whole (now-larger) date column so the existing numFmt is reused/updated.
$ws.Range("A2:A39").NumberFormat = "yyyy-mm-dd h:mm:ss"
